$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.509.04'
$ws.Range('E2').Value = '  +2.07%  '

$ws.Range('D3').Value = '2.527.85'
$ws.Range('E3').Value = '  +2.45%  '

$ws.Range('E4').Value = '  -0.18%  '

$ws.Range('D5').Value = '''594.16'
$ws.Range('E5').Value = '  +1.96%  '

$ws.Range('D6').Value = '''177.35'
$ws.Range('E6').Value = '  +1.89%  '

$ws.Range('E7').Value = '  -0.10%  '

$ws.Range('D8').Value = '''0.519'
$ws.Range('E8').Value = '  +1.30%  '

$ws.Range('D9').Value = '2.528.25'
$ws.Range('E9').Value = '  +2.48%  '

$ws.Range('D10').Value = '''0.144'
$ws.Range('E10').Value = '  +4.69%  '

$ws.Range('E11').Value = '  -1.03%  '

$ws.Range('D12').Value = '''5.00'
$ws.Range('E12').Value = '  +1.08%  '

$ws.Range('D13').Value = '''0.337'
$ws.Range('E13').Value = '  +1.45%  '

$ws.Range('D14').Value = '2.965.10'
$ws.Range('E14').Value = '  +1.24%  '

$ws.Range('D15').Value = '''25.94'
$ws.Range('E15').Value = '  +2.23%  '

$ws.Range('D16').Value = '68.437.40'
$ws.Range('E16').Value = '  +1.95%  '

$ws.Range('E17').Value = '  +0.82%  '

$ws.Range('D18').Value = '2.514.76'
$ws.Range('E18').Value = '  +2.33%  '

$ws.Range('D19').Value = '''11.07'
$ws.Range('E19').Value = '  +1.55%  '

$ws.Range('D20').Value = '''7.50'
$ws.Range('E20').Value = '  +0.56%  '

$ws.Range('D21').Value = '''351.99'
$ws.Range('E21').Value = '  +0.87%  '

$ws.Range('D22').Value = '''4.13'
$ws.Range('E22').Value = '  +2.64%  '

$ws.Range('D23').Value = '''71.26'
$ws.Range('E23').Value = '  +2.73%  '

$ws.Range('E24').Value = '  +0.03%  '

$ws.Range('D25').Value = '''4.22'
$ws.Range('E25').Value = '  +0.75%  '

$ws.Range('D26').Value = '''1.73'
$ws.Range('E26').Value = '  -3.65%  '

$ws.Range('D27').Value = '''9.16'
$ws.Range('E27').Value = '  +0.14%  '

$ws.Range('D28').Value = '2.640.40'
$ws.Range('E28').Value = '  +1.81%  '

$ws.Range('D29').Value = '''0.998'
$ws.Range('E29').Value = '  -0.01%  '

$ws.Range('D30').Value = '''514.73'
$ws.Range('E30').Value = '  +2.78%  '

$ws.Range('D31').Value = '0.0₃0901'
$ws.Range('E31').Value = '  -0.07%  '

$ws.Range('D32').Value = '''7.83'
$ws.Range('E32').Value = '  +1.29%  '

$ws.Range('D33').Value = '''1.27'
$ws.Range('E33').Value = '  +2.78%  '

$ws.Range('D34').Value = '''1.78'
$ws.Range('E34').Value = '  +1.46%  '

$ws.Range('D35').Value = '''0.999'
$ws.Range('E35').Value = '  -0.06%  '

$ws.Range('E36').Value = '  +1.41%  '

$ws.Range('D37').Value = '''162.13'
$ws.Range('E37').Value = '  +0.45%  '

$ws.Range('E38').Value = '  +0.06%  '

$ws.Range('D39').Value = '''18.36'
$ws.Range('E39').Value = '  +1.41%  '

$ws.Range('D40').Value = '''1.33'
$ws.Range('E40').Value = '  +0.37%  '

$ws.Range('E41').Value = '  +0.04%  '

$ws.Range('D42').Value = '''1.76'
$ws.Range('E42').Value = '  +4.21%  '

$ws.Range('D43').Value = '''0.329'
$ws.Range('E43').Value = '  +0.61%  '

$ws.Range('D44').Value = '''4.85'
$ws.Range('E44').Value = '  +0.56%  '

$ws.Range('D45').Value = '''2.44'
$ws.Range('E45').Value = '  +2.21%  '

$ws.Range('D46').Value = '''151.87'
$ws.Range('E46').Value = '  +6.87%  '

$ws.Range('D47').Value = '''3.57'
$ws.Range('E47').Value = '  +2.98%  '

$ws.Range('D48').Value = '''0.521'
$ws.Range('E48').Value = '  +2.13%  '

$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₆0260'
$ws.Range('E49').Value = '  +2.14%  '

$ws.Range('B50').Value = 'Optimism'
$ws.Range('C50').Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range('D50').Value = '''1.61'
$ws.Range('E50').Value = '  +2.43%  '

$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '''0.0742'
$ws.Range('E51').Value = '  +0.45%  '
